{"js": "// Apply the resume edits described by the commit:\n//   \"Added RTOS as skill, changed from ROS2 to ROS\"\n// plus the accompanying small wording tweaks that show up as genuine\n// (non-cosmetic) <w:t> text changes in the target OOXML.\n//\n// NOTE: the diff also contains a large amount of proofErr (spell-check\n// squiggle) removal / run-merging noise that does not alter any visible\n// text anywhere in the document (every proofErr-wrapped word in the\n// original file is touched). That is a side effect of the document being\n// re-saved, not an intentional edit, so it is intentionally NOT\n// reproduced here \u2014 only the rendered text changes are applied.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  // Only touch the first match - some search strings are intentionally\n  // specific enough to be unique, but guard anyway.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Skill line: \"ROS2,\" -> \"ROS,\". \"ROS2,\" also appears later in an\n// experience bullet (\"...using ROS2, CARLA and CI principles\"), which\n// must stay untouched, so anchor the search with enough left-hand\n// context to stay unique to the skills line.\nawait replaceOnce(\n  \"MATLAB), ROS2,\",\n  \"MATLAB), ROS,\"\n);\n\n// 2) Add \"RTOS, \" as a new skill right after \"software testing\n// (GoogleTest, pytest), \" and before \"CI/CD, \". Anchored on \"), CI/CD, \"\n// so the GoogleTest/pytest runs (and their proofErr spell-check marks)\n// are left completely untouched.\nawait replaceOnce(\n  \"), CI/CD, \",\n  \"), RTOS, CI/CD, \"\n);\n\n// 3) Remove the \"(oscilloscope, multimeter)\" parenthetical after \"debugging\".\nawait replaceOnce(\n  \"debugging (oscilloscope, multimeter) \",\n  \"debugging  \"\n);\n\n// 4) \"for Level 2 vehicle autonomy\" -> \"for level 2 vehicle autonomy\"\nawait replaceOnce(\n  \"for Level 2 vehicle autonomy\",\n  \"for level 2 vehicle autonomy\"\n);\n\n// 5) \"Working on structured testing framework\" -> \"Currently working on structured testing framework\"\nawait replaceOnce(\n  \"Working on structured testing framework\",\n  \"Currently working on structured testing framework\"\n);\n", "ps1": "# Apply the resume edits described by the commit:\n#   \"Added RTOS as skill, changed from ROS2 to ROS\"\n# plus the accompanying small wording tweaks that show up as genuine\n# (non-cosmetic) <w:t> text changes in the target OOXML.\n#\n# NOTE: the target diff also contains a large amount of proofErr\n# (spell-check squiggle) removal / run-merging noise that does not alter\n# any visible text anywhere in the document (every proofErr-wrapped word\n# in the original file is touched by it). That is a side effect of the\n# document being re-saved, not an intentional edit, so it is\n# intentionally NOT reproduced here - only the rendered text changes are\n# applied, using narrow Find/Replace anchors that avoid disturbing the\n# proofErr-wrapped runs.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue = 1, wdReplaceOne = 1 (we only ever want the single match)\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\n# 1) Skill line: \"ROS2,\" -> \"ROS,\". \"ROS2,\" also appears later in an\n# experience bullet (\"...using ROS2, CARLA and CI principles\"), which\n# must stay untouched, so anchor the search with enough left-hand\n# context to stay unique to the skills line.\nReplace-Text \"MATLAB), ROS2,\" \"MATLAB), ROS,\"\n\n# 2) Add \"RTOS, \" as a new skill right after \"software testing\n# (GoogleTest, pytest), \" and before \"CI/CD, \". Anchored on \"), CI/CD, \"\n# so the GoogleTest/pytest runs (and their proofErr spell-check marks)\n# are left completely untouched.\nReplace-Text \"), CI/CD, \" \"), RTOS, CI/CD, \"\n\n# 3) Remove the \"(oscilloscope, multimeter)\" parenthetical after \"debugging\".\nReplace-Text \"debugging (oscilloscope, multimeter) \" \"debugging  \"\n\n# 4) \"for Level 2 vehicle autonomy\" -> \"for level 2 vehicle autonomy\"\nReplace-Text \"for Level 2 vehicle autonomy\" \"for level 2 vehicle autonomy\"\n\n# 5) \"Working on structured testing framework\" -> \"Currently working on structured testing framework\"\nReplace-Text \"Working on structured testing framework\" \"Currently working on structured testing framework\"\n"}
